$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date text change
$ws.Range("B1").Value = "14/03/2023"

# Update B (values) and C (hour) columns for rows 2-20
$ws.Range("B2").Value = 525.2
$ws.Range("C2").Value = 38

$ws.Range("B3").Value = 236
$ws.Range("C3").Value = 38

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 38

$ws.Range("B5").Value = 96
$ws.Range("C5").Value = 38

$ws.Range("B6").Value = 38
$ws.Range("C6").Value = 38

$ws.Range("B7").Value = 82
$ws.Range("C7").Value = 38

$ws.Range("B8").Value = 35

$ws.Range("B9").Value = 77

$ws.Range("B10").Value = 440
$ws.Range("C10").Value = 38

$ws.Range("B11").Value = 148
$ws.Range("C11").Value = 38

$ws.Range("B12").Value = 483.6
$ws.Range("C12").Value = 38

$ws.Range("B13").Value = 355
$ws.Range("C13").Value = 38

$ws.Range("B14").Value = 530
$ws.Range("C14").Value = 38

$ws.Range("B15").Value = 159
$ws.Range("C15").Value = 38

$ws.Range("B16").Value = 124

$ws.Range("B17").Value = 86
$ws.Range("C17").Value = 38

$ws.Range("C18").Value = 38

$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 38

$ws.Range("B20").Value = 59
$ws.Range("C20").Value = 38
